# Append the latest Adafruit IO reading as row 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C:F ("24", "0.0", "0.0", "0.0") look numeric to Excel's
# auto-detection, so force them to Text before writing and then restore
# the cell style afterward (matches the source file, where every value -
# numeric-looking or not - is stored as literal text).
$ws.Range("C4:F4").NumberFormat = "@"

$ws.Range("A4").Value = "2024-09-25T17:55:14Z"
$ws.Range("B4").Value = "temperature"
$ws.Range("C4").Value = "24"
$ws.Range("D4").Value = "0.0"
$ws.Range("E4").Value = "0.0"
$ws.Range("F4").Value = "0.0"

# Drop the temporary Text format so the new cells end up styled like the
# rest of the sheet (no explicit style index).
$ws.Range("C4:F4").Style = "Normal"
